$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shared-string text fix: "NIGE 1.1/VFX" -> "NIGE 1.1 / VFX" ---
$ws.Range("G1").Value = "NIGE 1.1 / VFX"

# --- 2. Column G width 15.5703125 -> 18 (character units) ---
# COM ColumnWidth is offset from the stored OOXML width by 5/6, so subtract it.
$ws.Columns.Item(7).ColumnWidth = 18 - 5/6

# --- 3. Updated input (raw) values ---
$ws.Range("B17").Value = 1
$ws.Range("E13").Value = 7207
$ws.Range("E18").Value = 9410
$ws.Range("E19").Value = 35643
$ws.Range("E20").Value = 33387
$ws.Range("E21").Value = 26968

# --- 4. Flip the F-column ratio formulas from D/E to E/D ---
$ws.Range("F7").Formula = "=E7/D7"
$ws.Range("F8").Formula = "=E8/D8"
$ws.Range("F9").Formula = "=E9/D9"
$ws.Range("F10").Formula = "=E10/D10"
$ws.Range("F11").Formula = "=E11/D11"
$ws.Range("F12").Formula = "=E12/D12"
$ws.Range("F13").Formula = "=E13/D13"
$ws.Range("F14").Formula = "=E14/D14"

$ws.Range("F16").Formula = "=E16/D16"
$ws.Range("F17").Formula = "=E17/D17"
$ws.Range("F18").Formula = "=E18/D18"
$ws.Range("F19").Formula = "=E19/D19"
$ws.Range("F20").Formula = "=E20/D20"
$ws.Range("F21").Formula = "=E21/D21"
$ws.Range("F22").Formula = "=E22/D22"

# --- 5. The "bordered percent" look (style used for the last line item of each
#        block) now belongs to the SUM row beneath it, since that row joined the
#        shared formula group. Swap the cell formats between the old/new
#        border-bearing rows, reusing the existing style slots.
$ws.Range("F13").Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F13").PasteSpecial(-4122) | Out-Null

$ws.Range("F21").Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$ws.Range("F17").Copy() | Out-Null
$ws.Range("F21").PasteSpecial(-4122) | Out-Null

# --- 6. Move the active selection to B26 ---
[void]$ws.Range("B26").Select()
